$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 844.3
$ws.Range("J28").Value = 991.3333
$ws.Range("L28").Value = 991.3333
$ws.Range("N28").Value = -1961.3333

# ALC row 39
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 278.9
$ws.Range("I39").Value = 180.14285
$ws.Range("J39").Value = 509.33334
$ws.Range("K39").Value = 540.4285500000001
$ws.Range("L39").Value = 1528.00002
$ws.Range("M39").Value = -244.4285500000001
$ws.Range("N39").Value = -2120.00002

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1035
$ws.Range("J40").Value = 984.3333
$ws.Range("L40").Value = 984.3333
$ws.Range("N40").Value = -1334.3333

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 461.69232
$ws.Range("I80").Value = 379.8
$ws.Range("J80").Value = 512.875
$ws.Range("K80").Value = 1139.4
$ws.Range("L80").Value = 1538.625
$ws.Range("M80").Value = -141.4000000000001
$ws.Range("N80").Value = -3534.625

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 461.69232
$ws.Range("I83").Value = 379.8
$ws.Range("J83").Value = 512.875
$ws.Range("K83").Value = 3418.2
$ws.Range("L83").Value = 4615.875
$ws.Range("M83").Value = 1573.8
$ws.Range("N83").Value = -14599.875

# ARM row 63
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1758
$ws.Range("I63").Value = 1447.5
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1447.5
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -761.5
$ws.Range("N63").Value = -4372

# ARM row 66
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 1758
$ws.Range("I66").Value = 1447.5
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 7237.5
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -3805.5
$ws.Range("N66").Value = -21864

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1755.5
$ws.Range("J64").Value = 1674
$ws.Range("L64").Value = 1674
$ws.Range("N64").Value = -2124

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 1755.5
$ws.Range("J67").Value = 1674
$ws.Range("L67").Value = 1674
$ws.Range("N67").Value = -3234

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2685.1667
$ws.Range("I134").Value = 822.2
$ws.Range("K134").Value = 2466.6
$ws.Range("M134").Value = 68.39999999999964

# CRP row 7
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 32.2
$ws.Range("I7").Value = 18
$ws.Range("K7").Value = 18
$ws.Range("M7").Value = 95

# CRP row 47
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 19300
$ws.Range("I47").Value = 19300
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 19300
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -18734
$ws.Range("N47").ClearContents()

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 3385.5715
$ws.Range("J94").Value = 924.5
$ws.Range("L94").Value = 924.5
$ws.Range("N94").Value = -1826.5

# CUL row 44
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 869.25
$ws.Range("I44").Value = 294.5
$ws.Range("J44").Value = 1444
$ws.Range("K44").Value = 883.5
$ws.Range("L44").Value = 4332
$ws.Range("M44").Value = -485.5
$ws.Range("N44").Value = -5128

# CUL row 92
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 494.5
$ws.Range("J92").Value = 494.5
$ws.Range("L92").Value = 1483.5
$ws.Range("N92").Value = -3979.5

# CUL row 95
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H95").Value = 12500
$ws.Range("J95").Value = 12500
$ws.Range("L95").Value = 37500
$ws.Range("N95").Value = -41618

# CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 441.66666
$ws.Range("I103").Value = 441.66666
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 1324.99998
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -445.9999800000001
$ws.Range("N103").ClearContents()

# CUL row 108
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 317.33334
$ws.Range("I108").Value = 317.33334
$ws.Range("K108").Value = 952.0000200000001
$ws.Range("M108").Value = 1927.99998

# CUL row 117
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 5467.4287
$ws.Range("I117").Value = 642
$ws.Range("J117").Value = 9086.5
$ws.Range("K117").Value = 1926
$ws.Range("L117").Value = 27259.5
$ws.Range("M117").Value = 1516
$ws.Range("N117").Value = -34143.5

# CUL row 118
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 1900.6666
$ws.Range("I118").Value = 2074.8
$ws.Range("K118").Value = 6224.400000000001
$ws.Range("M118").Value = -4981.400000000001

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 15000
$ws.Range("J125").Value = 20000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

# CUL row 126
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 1500
$ws.Range("I126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = 440

# CUL row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 5220.1665
$ws.Range("I139").Value = 2257.6
$ws.Range("K139").Value = 6772.799999999999
$ws.Range("M139").Value = -1632.799999999999

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 367.83334
$ws.Range("I140").Value = 367.83334
$ws.Range("K140").Value = 1103.50002
$ws.Range("M140").Value = 4076.49998

# GSM row 24
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 17172.273
$ws.Range("J80").Value = 22470
$ws.Range("L80").Value = 22470
$ws.Range("N80").Value = -24466

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 17172.273
$ws.Range("J83").Value = 22470
$ws.Range("L83").Value = 112350
$ws.Range("N83").Value = -122334

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1197.091
$ws.Range("I97").Value = 1070.7142
$ws.Range("J97").Value = 1418.25
$ws.Range("K97").Value = 1070.7142
$ws.Range("L97").Value = 1418.25
$ws.Range("M97").Value = -574.7141999999999
$ws.Range("N97").Value = -2410.25

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 599.5
$ws.Range("J82").Value = 599.5
$ws.Range("L82").Value = 599.5
$ws.Range("N82").Value = -1321.5

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 599.5
$ws.Range("J85").Value = 599.5
$ws.Range("L85").Value = 599.5
$ws.Range("N85").Value = -3095.5

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 820
$ws.Range("I107").Value = 892
$ws.Range("K107").Value = 2676
$ws.Range("M107").Value = -756

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2499.4167
$ws.Range("I136").Value = 2149.3
$ws.Range("J136").Value = 4250
$ws.Range("K136").Value = 6447.900000000001
$ws.Range("L136").Value = 12750
$ws.Range("M136").Value = -3897.900000000001
$ws.Range("N136").Value = -17850
